# language.xlsx — "some tweaks to level 2-2 ... text tweaks"
#
# Updates the wording of several localized UI/dialogue strings on the "en"
# sheet (column B holds the text for the key in column A), and moves the
# active selection/scroll position to reflect where the author was last
# working in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# level01_intro_2
$ws.Range("B59").Value = "Each object has properties that fit different jobs."
# level01_classify_2
$ws.Range("B64").Value = "This allows us to use the right object for every job!"
# level01_play_begin_1
$ws.Range("B69").Value = "Look! There's a star piece."
# level02_heavy_1 (also fixes missing space: "heavycompared" -> "heavy compared")
$ws.Range("B75").Value = "A 40-pound iron block is heavy compared to a 20-pound block."
# level02_begin_1 (added comma after "Now")
$ws.Range("B76").Value = "Now, let's sort objects into two classifications: light and heavy!"
# level03_intro_1 (added comma after "Now")
$ws.Range("B79").Value = "Now, let's look at buoyancy."
# level03_buoyancy_2 (rewritten)
$ws.Range("B81").Value = "If an object has a lot of buoyancy, it will move up and float above water."
# level03_buoyancy_3 (rewritten)
$ws.Range("B82").Value = "If an object has no buoyancy, it will move down and sink under water."
# level04_intro_1
$ws.Range("B84").Value = "The final mission! We will learn about the conductive property."
# level04_non_conductive_1 (rewritten)
$ws.Range("B87").Value = "The granite rock, on the other hand, is not conductive. The electricity is unable to flow through at all."

# Move the sheet's scroll/selection to where the edits were made (B82).
$win = $excel.ActiveWindow
$win.ScrollRow = 70
$win.ScrollColumn = 1
[void]$ws.Range("B82").Select()
